# Shift the quarterly data one period to the left (drop the oldest
# period, which lived in column D) and append the newly published
# period ("12 ماهه منتهی به 1401/12", published 1402-02-30) in column M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Delete column D entirely - this shifts E:M left to D:L, carrying
# values/styles/number formats/column widths with them.
$ws.Columns("D").Delete()

# New rightmost column is now M; give it the "annual period" width (29)
# like the other 12-month columns (F, J before the shift -> now E, I).
$ws.Columns("M").ColumnWidth = 29

# Period header (row 8) and publish-date header (row 9) for the new column.
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-30"

# New financial figures for the just-published period.
$ws.Range("M11").Value = 14721432
$ws.Range("M12").Value = -13158971
$ws.Range("M13").Value = 1562461
$ws.Range("M14").Value = -265630
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = -9
$ws.Range("M17").Value = 1296822
$ws.Range("M18").Value = -34926
$ws.Range("M19").Value = 242640
$ws.Range("M20").Value = 1504536
$ws.Range("M21").Value = -224314
$ws.Range("M22").Value = 1280222
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 1280222
$ws.Range("M25").Value = 6577

# Row 26 ("سرمایه") was recomputed with an updated read_price algorithm,
# so it is not a pure left-shift: fix up the I26 figure that the column
# delete brought over from the old J26, then set the freshly computed M26.
$ws.Range("I26").Value = 194650
$ws.Range("M26").Value = 194650

$ws.Range("M27").Value = 6565
